# One-click update from Van Paper 07:10 AM on 2026-01-02
#
# - Renames customer "NAI LEGACY" (row 38) to "TONKA PARTNERS LLC"
# - Inserts two new prospect rows (1655 BEAM LLC / LIQUOR BARREL CP)
#   right after the TONKA PARTNERS LLC row, pushing the remaining
#   leaderboard rows down by two (old rows 39-42 become 41-44)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing customer on row 38
$ws.Range("A38").Value = "TONKA PARTNERS LLC"

# Insert two fresh rows right below row 38 (before the old row 39)
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# Match the row height used throughout the rest of the leaderboard
$ws.Rows.Item(39).RowHeight = 13.05
$ws.Rows.Item(40).RowHeight = 13.05

# New row 39: 1655 BEAM LLC
$ws.Range("A39").Value = "1655 BEAM LLC"
$ws.Range("B39").Value = "Ballman, John W"
$ws.Range("C39").Value = "015"
$ws.Range("E39").Value = "0008399"

# New row 40: LIQUOR BARREL CP
$ws.Range("A40").Value = "LIQUOR BARREL CP"
$ws.Range("B40").Value = "Larsen, Rick J"
$ws.Range("C40").Value = "040"
$ws.Range("E40").Value = "0008400"
